# Apply updated crypto price/volume figures (Sheet1, columns D & E, rows 2-51)
# to match the scraped commit's data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.177.68"
$ws.Range("E2").Value = "  -1.95%  "
$ws.Range("D3").Value = "3.589.59"
$ws.Range("E3").Value = "  -2.74%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("E5").Value = "  +20.51%  "
$ws.Range("D6").Value = "'224.48"
$ws.Range("E6").Value = "  -5.17%  "
$ws.Range("D7").Value = "'633.35"
$ws.Range("E8").Value = "  -3.59%  "
$ws.Range("D9").Value = "'1.07"
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("D11").Value = "3.589.95"
$ws.Range("E11").Value = "  -2.67%  "
$ws.Range("D12").Value = "'46.37"
$ws.Range("E12").Value = "  +4.95%  "
$ws.Range("E13").Value = "  -1.37%  "
$ws.Range("D14").Value = "'0.0000286"
$ws.Range("E14").Value = "  -5.82%  "
$ws.Range("E15").Value = "  -4.99%  "
$ws.Range("D16").Value = "4.260.87"
$ws.Range("E16").Value = "  -2.71%  "
$ws.Range("D17").Value = "94.901.87"
$ws.Range("E17").Value = "  -1.93%  "
$ws.Range("E18").Value = "  -4.14%  "
$ws.Range("D19").Value = "'19.65"
$ws.Range("E19").Value = "  +4.92%  "
$ws.Range("D20").Value = "3.592.86"
$ws.Range("E20").Value = "  -2.68%  "
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("D23").Value = "'498.37"
$ws.Range("E23").Value = "  -4.13%  "
$ws.Range("E24").Value = "  -6.09%  "
$ws.Range("D25").Value = "'0.239"
$ws.Range("E25").Value = "  +19.60%  "
$ws.Range("D26").Value = "'117.23"
$ws.Range("E26").Value = "  +15.62%  "
$ws.Range("E27").Value = "  -4.41%  "
$ws.Range("E28").Value = "  -3.26%  "
$ws.Range("D29").Value = "3.781.39"
$ws.Range("E29").Value = "  -2.82%  "
$ws.Range("D30").Value = "'12.46"
$ws.Range("E30").Value = "  -7.15%  "
$ws.Range("D31").Value = "'12.95"
$ws.Range("E31").Value = "  +3.31%  "
$ws.Range("D32").Value = "'2.88"
$ws.Range("E32").Value = "  -4.68%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D36").Value = "'1.74"
$ws.Range("E36").Value = "  -6.92%  "
$ws.Range("D37").Value = "'31.44"
$ws.Range("E37").Value = "  -2.31%  "
$ws.Range("D38").Value = "'0.581"
$ws.Range("E38").Value = "  -1.76%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").Value = "'587.39"
$ws.Range("E40").Value = "  -9.36%  "
$ws.Range("E41").Value = "  -6.62%  "
$ws.Range("D42").Value = "'6.78"
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("D43").Value = "'40.77"
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("E44").Value = "  -2.39%  "
$ws.Range("E45").Value = "  -5.64%  "
$ws.Range("D46").Value = "'0.0467"
$ws.Range("E46").Value = "  +0.64%  "
$ws.Range("D47").Value = "'1.89"
$ws.Range("E47").Value = "  -7.72%  "
$ws.Range("D48").Value = "'0.911"
$ws.Range("E48").Value = "  -5.14%  "
$ws.Range("D49").Value = "'23.47"
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("D50").Value = "'3.64"
$ws.Range("E50").Value = "  +2.91%  "
$ws.Range("D51").Value = "'221.53"
$ws.Range("E51").Value = "  +8.19%  "
